$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9570517539978027
$ws.Range("B1").Value = 1.661927342414856
$ws.Range("C1").Value = 3.302589416503906
$ws.Range("D1").Value = 2.619236946105957
$ws.Range("E1").Value = 0.3567759394645691
